# ---------------------------------------------------------------------------
# Applies two changes to "My Skills.docx":
#   1. Collapse the "I am VERY detail oriented..." paragraph's three runs
#      (which were split around a spell-check-flagged "so.") into a single
#      plain run with the full sentence, dropping the w:proofErr markers.
#   2. Append a page break plus a new "Frameworks/tools" section (with
#      FE / BE / other-tools sub-lists) at the end of the document, after
#      the existing last paragraph ("... to reflect my current skillset
#      and value adds"), leaving that paragraph's runs untouched.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1: merge the "detail oriented" runs into one run -----------------
# The original sentence spans three runs/proofErr-wrapped fragments that,
# concatenated, read "I am VERY detail oriented. Sometimes maybe too much
# so. ". Searching for that full concatenation lets Find match across the
# run/proofErr boundaries and replace the whole span with one plain run.
$d.Content.Find.Execute(
    "I am VERY detail oriented. Sometimes maybe too much so. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I am VERY detail oriented. Sometimes maybe too much so. ", 2) | Out-Null

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Add-Paragraph-Xml([string]$innerXml) {
    # Inserts a brand-new <w:p> at the very end of the document body
    # (just before the sectPr) without disturbing any existing content.
    $e = $d.Content.End
    $r = $d.Range($e, $e)
    $xml = "<w:p $ns>$innerXml</w:p>"
    $r.InsertXML($xml) | Out-Null
}

# --- Change 2: append the new trailing content --------------------------------

# Page break paragraph
Add-Paragraph-Xml "<w:r><w:br w:type='page'/></w:r>"

# "Frameworks/tools" heading paragraph (carries the rendered-page-break marker)
Add-Paragraph-Xml "<w:r><w:lastRenderedPageBreak/><w:t>Frameworks/tools</w:t></w:r>"

# FE
Add-Paragraph-Xml "<w:r><w:t>FE</w:t></w:r>"

Add-Paragraph-Xml (
    "<w:r><w:t xml:space='preserve'>react, next.js, typescript, </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>javascript</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'>, </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>vue</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'>, tailwind </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>css</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>"
)

# BE
Add-Paragraph-Xml "<w:r><w:t>BE</w:t></w:r>"

Add-Paragraph-Xml (
    "<w:r><w:t xml:space='preserve'>node, express, </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>postgresql</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'>, rest </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>api</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t>, Scaffold eth, python</w:t></w:r>"
)

# other tools
Add-Paragraph-Xml "<w:r><w:t>other tools</w:t></w:r>"

Add-Paragraph-Xml (
    "<w:r><w:t xml:space='preserve'>git, postman, docker, </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>graphql</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'>, windsurf, </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>npm</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>"
)
